# "refractored partial year into separate method"
#
# B3 used to bake the partial-year factor (C27) straight into the gross
# "Navy Ret" income formula (70*1500*C27) and then divide it back out again
# in C3 (B3/C27) to recover the full-year figure. Pull the partial-year
# factor out of B3/C3 and apply it later where the partial year actually
# belongs (B12 and B21), and re-tune the tax bracket constants used in B13
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- relabel the two rows whose headings were tweaked -----------------
# (set A18 before A15 so new shared-string entries land in the same order
# as the target: "Gross Exp " then "Expenses (final_annual_expenses)")
$ws.Range("A18").Value = "Gross Exp "
$ws.Range("A15").Value = "Expenses (final_annual_expenses)"

# --- pull C27 (partial-year fraction) out of B3 / C3 -------------------
$ws.Range("B3").Formula = "=70*1500"
$ws.Range("C3").Formula = "=B3"

# --- apply the partial-year fraction down in B12 / C12 instead ---------
$ws.Range("B12").Formula = "=(B10-B20)*C27"
$ws.Range("C12").Formula = "=C10-C20"

# --- re-tune the tax-bracket constants used for the first (partial) year
$ws.Range("B13").Formula = "=ROUND((1700 + 0.15*(B12- 36000) + B3 * 0.0765),0)"

# --- B21 also needs the partial-year fraction applied to B10 -----------
$ws.Range("B21").Formula = "=B10*C27-B13-B18"

# --- yearly IRA->non-IRA transfer amount dropped slightly --------------
$ws.Range("B26").Value = 14500

# --- best-effort reproduction of the refreshed window/pane scroll state
$win = $excel.ActiveWindow
$win.SplitColumn = 1
$win.SplitRow = 1
$win.FreezePanes = $true

$win.Panes.Item(1).Activate()
$ws.Range("A19").Select()

$win.Panes.Item(2).Activate()
$ws.Range("A2").Select()

$win.Panes.Item(3).Activate()
$ws.Range("A15").Select()

$win.Panes.Item(4).Activate()
$ws.Range("L24").Select()

$wb.Application.Calculate()
